$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to move from
# 45179 (2023-09-10) to 45180 (2023-09-11) for every data row (2-439).
$ws.Range("C2:C439").Value = 45180
